# Update automatico via Actualizar 02-04-2021 18-35-09
#
# This mirrors the automated "availability check" refresh: a brand-new
# timestamp snapshot is recorded at the top of the history block (rows
# 2-15), and the previous snapshots cascade one block down (old rows
# 2-15 -> rows 16-29, old rows 16-29 -> rows 30-37). Only column D
# (Fecha) changes; everything else is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New snapshot timestamp written into the first block (rows 2-15).
$newTimestamp = 44231.77435050038

# Value that the first block (rows 2-15) held before this update -
# it now becomes the value for the second block (rows 16-29).
$shiftedBlock1 = 44231.76464553241

# Value that the second block (rows 16-29) held before this update -
# it now becomes the value for the third block (rows 30-37).
$shiftedBlock2 = 44231.76190935185

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $newTimestamp
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $shiftedBlock1
}

for ($r = 30; $r -le 37; $r++) {
    $ws.Cells.Item($r, 4).Value = $shiftedBlock2
}
